# Update "想去人数" (interest count) figures in column F.
# The "展览" sheet and the "全部类型" sheet both list the same exhibition
# rows (the latter also includes one extra show row from "演出", which
# shifts everything below it down by one row), so each sheet gets its
# own row -> new-value map.

$wb = $excel.ActiveWorkbook

# "展览" sheet (row -> new F value)
$updatesExhibition = @{
    2  = 203
    3  = 186
    4  = 5265
    8  = 593
    9  = 554
    13 = 4208
    15 = 177
    18 = 3292
    19 = 160
    20 = 1074
    23 = 191
    25 = 34
    26 = 137
    28 = 298
    29 = 27
    32 = 20
    33 = 22
}

# "全部类型" sheet (row -> new F value) - same rows, shifted down by 1
$updatesAll = @{
    2  = 203
    3  = 186
    5  = 5265
    9  = 593
    10 = 554
    14 = 4208
    16 = 177
    19 = 3292
    20 = 160
    21 = 1074
    24 = 191
    26 = 34
    27 = 137
    29 = 298
    30 = 27
    33 = 20
    34 = 22
}

$ws = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibition.Keys) {
    $ws.Range("F$row").Value = $updatesExhibition[$row]
}

$ws = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAll.Keys) {
    $ws.Range("F$row").Value = $updatesAll[$row]
}
